# Apply the "updated 4.0 files and mdl" change:
#  - About sheet: update the date in C1 from 2024-01-29 to 2024-04-10
#  - MCF sheet: update the capacity-factor values in column B (rows 2-18, except
#    rows that are already 1 or 0) to 1, which cascades into the dependent
#    formulas in B19:B25
#  - MCF sheet: update the selected/active cell to B17

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsMCF = $wb.Worksheets.Item("MCF")

# Update the "About" sheet date stamp in C1 (serial 45320 -> 45392, i.e. 2024-04-10)
$wsAbout.Range("C1").Value = Get-Date -Year 2024 -Month 4 -Day 10 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

# Update the Maximum Capacity Factor values on the MCF sheet to 1
$rows = @(2, 3, 4, 6, 10, 11, 12, 13, 14, 16, 17, 18)
foreach ($r in $rows) {
    $wsMCF.Cells.Item($r, 2).Value = 1
}

# Move the active selection on the MCF sheet to B17
$wsMCF.Activate()
$wsMCF.Range("B17").Select()
